# Insert a new data row at row 541 (pushing existing rows 541-639 down to 542-640)
# and populate it with a new "Ajo" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(541).Insert()

$ws.Range("A541").Value = 10
$ws.Range("B541").Value = "Vega Modelo de Temuco"
$ws.Range("C541").Value = "La Araucanía"
$ws.Range("D541").Value = 44694
$ws.Range("D541").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E541").Value = 9
$ws.Range("F541").Value = 100112003
$ws.Range("G541").Value = "Ajo"
$ws.Range("H541").Value = "Chino"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 210
$ws.Range("K541").Value = 20000
$ws.Range("L541").Value = 22000
$ws.Range("M541").Value = 21190
$ws.Range("N541").Value = "`$/caja 10 kilos"
$ws.Range("O541").Value = "China"
$ws.Range("P541").Value = 2119
$ws.Range("Q541").Value = 10
$ws.Range("R541").Value = "Hortaliza"
